# Current Tasks Backlog - flesh out the "Model" block with a Performer
# column + 3 new tasks, add a new "Презентер" section, and push the
# "View" / "Data Base" / "Client-Server" sections further down the
# sheet as their own highlighted blocks. Finally resize the backing
# table and nudge the view down to where the new content lives.
#
# NOTE: the order in which the text cells below are written matters -
# it drives the order new entries land in xl/sharedStrings.xml - so
# values are intentionally poked in a specific sequence rather than
# strictly top-to-bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the old rows 3-16 so we can rebuild the backlog from a
# clean slate (the rows below the "Model" header are being renumbered).
$ws.Range("A3:C16").Clear()

# --- "Model" block (rows 3-6 keep their task names, now get a
#     Performer too) -------------------------------------------------
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Oriented Graph class realization"
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Block classes realization"
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Code to flowchart functionality (pizdos)"

$ws.Range("C3").Value = "Кочетов"
$ws.Range("C4").Value = "Чугунов"
$ws.Range("C5").Value = "Чугунов"

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Flowchart to code functionality (pizdos)"

# --- new "Презентер" section header (row 23), styled like the
#     "Model" header (row 2) ------------------------------------------
$ws.Range("A2:C2").Copy()
$ws.Range("A23:C23").PasteSpecial(-4122)

$ws.Range("B29").Value = "Взаимодействие с презентером"
$ws.Range("B23").Value = "Презентер"
$ws.Range("C7").Value = "Веселов"

$ws.Range("B24").Value = "Взаимодействие со вьюхами"
$ws.Range("C24").Value = "Веселов"

# --- rest of the "Model" block's new tasks (rows 7-11) ----------------
$ws.Range("B7").Value = "Построение модели"
$ws.Range("B8").Value = "Отрисовка всего"
$ws.Range("B9").Value = "Поддержка моделью координат"
$ws.Range("B10").Value = "Класс управление БД"
$ws.Range("B11").Value = "Реакция вьюхи на мелкие вещи ()"

# --- remaining "Презентер" section cells -------------------------------
$ws.Range("C29").Value = "Веселов/Кочетов"
$ws.Range("B25").Value = "Взаимодействие с моделью"
$ws.Range("C25").Value = "Веселов"
$ws.Range("C28").Value = "Кочетов"

# --- remaining numbers for the "Model" block ---------------------------
$ws.Range("A7").Value = 10
$ws.Range("A8").Value = 8
$ws.Range("A9").Value = 5
$ws.Range("A10").Value = 6
$ws.Range("A11").Value = 2

$ws.Range("A23").Value = 5

# --- "View" section header (row 27) + its tasks ------------------------
$ws.Range("A2:C2").Copy()
$ws.Range("A27:C27").PasteSpecial(-4122)
$ws.Range("A27").Value = 5
$ws.Range("B27").Value = "View"

$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Construct forms"

# --- "Data Base" section header (row 33) --------------------------------
$ws.Range("A2:C2").Copy()
$ws.Range("A33:C33").PasteSpecial(-4122)
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = "Data Base"

# --- "Client-Server" section header (row 35) ------------------------------
$ws.Range("A2:C2").Copy()
$ws.Range("A35:C35").PasteSpecial(-4122)
$ws.Range("A35").Value = 8
$ws.Range("B35").Value = "Client-Server"

$excel.CutCopyMode = $false

# Resize the backing table/ListObject to cover the new range.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C35"))

# Scroll the view down to the new content and move the selection.
$ws.Range("E22").Select()
